$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 18: new timesheet entry (Prabha, 02.01.2018) ---------------------
# C18 holds a dd.mm.yyyy-looking string ("02.01.2018"). Writing that text
# straight into a cell makes Excel auto-convert it to a date serial, so we
# add a trailing space first (defeats the date parser), strip the space via
# a TRIM() helper formula pasted back in as a value, then clean up the
# helper cell. Do this first so the new shared string for "02.01.2018" is
# registered before the other new string on this row.
$ws.Range("C18").Value = "02.01.2018 "
$ws.Range("Z1").Formula = "=TRIM(C18)"
$ws.Range("Z1").Copy()
$ws.Range("C18").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()

# Plain (non date-like) values can be written directly.
$ws.Range("A18").Value = 13
$ws.Range("B18").Value = "Prabha"
$ws.Range("D18").Value = "8.30 to 4.45"
$ws.Range("E18").Value = "project structural flow,process flow,api documentation"
$ws.Range("F18").Value = "completed"

# Match the formatting of the row above (font/fill/alignment) instead of
# whatever default formatting the new cells picked up.
$ws.Range("A17:E17").Copy()
$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("F17").Copy()
$ws.Range("F18").PasteSpecial(-4122)

$ws.Range("F20").Select() | Out-Null
